$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Mark chapter 3 (row 4) as "En proceso" for Walter (column B)
$ws.Range("B4").Value = "En proceso"

# Mark chapter 4 (row 5) as "Leido" for Walter (column B) - finalizacion de unidad 4
$ws.Range("B5").Value = "Leido"

# Update the active selection to C8
$ws.Range("C8").Select()
